$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.14

$ws.Range("B3").Value = 9.859999999999999
$ws.Range("D3").Value = 10.2

$ws.Range("C4").Value = 9.800000000000001
$ws.Range("E4").Value = 10.61

$ws.Range("D5").Value = 9.390000000000001
$ws.Range("G5").Value = 9.779999999999999

$ws.Range("G6").Value = 10.36
$ws.Range("H6").Value = 10.5

$ws.Range("E7").Value = 10.22
$ws.Range("F7").Value = 9.640000000000001
$ws.Range("I7").Value = 7.74

$ws.Range("F8").Value = 9.5

$ws.Range("G9").Value = 12.26
